# Updates the cryptos list price (D) and 1h volume change (E) columns
# for rows 2-51 on Sheet1, matching the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.121.80"
$ws.Range("E2").Value = "  +0.17%  "

$ws.Range("D3").Value = "2.967.49"
$ws.Range("E3").Value = "  -1.69%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'574.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.81%  "

$ws.Range("D6").Value = "'167.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.43%  "

$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "

$ws.Range("D8").Value = "'0.519"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").Value = "2.964.66"
$ws.Range("E9").Value = "  -1.67%  "

$ws.Range("D10").Value = "'6.70"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.19%  "

$ws.Range("D11").Value = "'0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.53%  "

$ws.Range("D12").Value = "'0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.24%  "

$ws.Range("D13").Value = "'0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.67%  "

$ws.Range("D14").Value = "'35.99"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.56%  "

$ws.Range("E15").Value = "  -0.46%  "

$ws.Range("D16").Value = "66.043.14"
$ws.Range("E16").Value = "  +0.14%  "

$ws.Range("D17").Value = "3.453.56"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").Value = "'7.15"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.77%  "

$ws.Range("D19").Value = "'16.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +14.87%  "

$ws.Range("D20").Value = "2.961.23"
$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("D21").Value = "'451.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.45%  "

$ws.Range("D22").Value = "'0.701"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.62%  "

$ws.Range("D23").Value = "'7.36"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.62%  "

$ws.Range("D24").Value = "'82.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.24%  "

$ws.Range("D25").Value = "'2.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.57%  "

$ws.Range("D26").Value = "'12.45"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("D27").Value = "'10.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.98%  "

$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").Value = "'8.34"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.03%  "

$ws.Range("D30").Value = "'2.49"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.71%  "

$ws.Range("D31").Value = "'2.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("D33").Value = "'0.117"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.88%  "

$ws.Range("D34").Value = "'27.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.43%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.16%  "

$ws.Range("D36").Value = "'0.981"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.31%  "

$ws.Range("D37").Value = "'5.77"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.35%  "

$ws.Range("D38").Value = "'47.55"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.51%  "

$ws.Range("D39").Value = "'49.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.37%  "

$ws.Range("D40").Value = "'2.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.15%  "

$ws.Range("D41").Value = "'0.306"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.79%  "

$ws.Range("D42").Value = "'0.122"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("D43").Value = "'2.87"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.15%  "

$ws.Range("D44").Value = "'8.56"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.24%  "

$ws.Range("D45").Value = "'385.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.35%  "

$ws.Range("D46").Value = "'0.0355"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "

$ws.Range("D47").Value = "2.691.51"
$ws.Range("E47").Value = "  -3.93%  "

$ws.Range("D48").Value = "'133.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.70%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("D50").Value = "'24.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("D51").Value = "'2.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.86%  "
